# Removed Test Case Inter-Dependency
#
# The product name / short name values get a distinguishing suffix so this
# sheet no longer collides with another test case, and the active tab
# switches from the input sheet ("ProductLoanInput") to the output sheet
# ("ProductLoanOutput") so the workbook opens on the output tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoanInput
$ws2 = $wb.Worksheets.Item(2)   # ProductLoanOutput

# Update the product name shown on both the input and output sheets so the
# two stay in sync (they both reference the same underlying text).
$ws1.Range("B1").Value = "2465-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-NONE-1st"
$ws2.Range("B1").Value = "2465-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-NONE-1st"

# Give the short name a unique value instead of the plain numeric code.
$ws1.Range("B2").Value = "2465d"

# Leave the input sheet selection resting on B1, and make the output sheet
# the active / selected tab of the workbook.
$ws1.Range("B1").Select() | Out-Null
$ws2.Activate()
